$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new data row right before the current row 164 ("Crespo
# record" / "Primera" entry dated 2021-12-02). Everything that used to live
# on rows 164..274 shifts down by one (to 165..275), exactly matching the
# target diff ("A1:R274" -> "A1:R275").
$ws.Rows.Item(164).Insert()

# Populate the newly inserted row 164 with the new weekly record.
$ws.Cells.Item(164, 1).Value = 7
$ws.Cells.Item(164, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(164, 3).Value = "Ñuble"
$ws.Cells.Item(164, 4).Value = 44824
$ws.Cells.Item(164, 5).Value = 16
$ws.Cells.Item(164, 6).Value = 100112006
$ws.Cells.Item(164, 7).Value = "Repollo"
$ws.Cells.Item(164, 8).Value = "Crespo record"
$ws.Cells.Item(164, 9).Value = "Primera"
$ws.Cells.Item(164, 10).Value = 120
$ws.Cells.Item(164, 11).Value = 1500
$ws.Cells.Item(164, 12).Value = 1600
$ws.Cells.Item(164, 13).Value = 1550
$ws.Cells.Item(164, 14).Value = "$/unidad"
$ws.Cells.Item(164, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(164, 16).Value = 1550
$ws.Cells.Item(164, 17).Value = 1
$ws.Cells.Item(164, 18).Value = "Hortaliza"

# Match the date-formatted style used by the rest of column D.
$ws.Cells.Item(164, 4).NumberFormat = $ws.Cells.Item(165, 4).NumberFormat
